$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: fill in the new timeline entry
$ws.Range("A16").Value = "Создание скрипта для копирования сотрудников из одной компании в другую с привязкой е ЛПУ и Регионам"
$ws.Range("B16").Value = 0.5

# Copy the date-formatted style from C15 onto C16 so it gets the same number format
$ws.Range("C15").Copy() | Out-Null
$ws.Range("C16").PasteSpecial(-4122) | Out-Null
$ws.Range("C16").Value = 43522

# Move the active selection to A17 (next empty row), matching the saved view state
$ws.Range("A17").Select() | Out-Null
